$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (AMM)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 878
$ws.Range("D2").Value = 910
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 1501.5
$ws.Range("J2").Value = -39.39393939393939

# Row 3 (IPR)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 84
$ws.Range("D3").Value = 89
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("I3").Value = 148
$ws.Range("J3").Value = -39.86486486486487

# Row 4 (MIG)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 213
$ws.Range("D4").Value = 226
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 2
$ws.Range("J4").Value = 88.33333333333333

# Row 5 (MOB)
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 518
$ws.Range("D5").Value = 530
$ws.Range("E5").Value = 1
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 926
$ws.Range("J5").Value = -42.7645788336933

# Row 6 (MOB PRE)
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 373
$ws.Range("D6").Value = 420
$ws.Range("E6").Value = 37
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 8
$ws.Range("I6").Value = 642
$ws.Range("J6").Value = -34.57943925233645

# Row 7 (MSK)
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 125
$ws.Range("D7").Value = 126
$ws.Range("I7").Value = 196
$ws.Range("J7").Value = -35.71428571428571

# Row 9 (TEC)
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 246
$ws.Range("D9").Value = 247
$ws.Range("G9").Value = 7
$ws.Range("I9").Value = 788
$ws.Range("J9").Value = -68.65482233502537

# Row 10 (TST)
$ws.Range("C10").Value = 48
$ws.Range("D10").Value = 48
$ws.Range("I10").Value = 95
$ws.Range("J10").Value = -49.47368421052632

# Row 11 (VIP)
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = -75

# Row 12 (WLC)
$ws.Range("C12").Value = 28
$ws.Range("D12").Value = 28
$ws.Range("I12").Value = 85
$ws.Range("J12").Value = -67.05882352941177
